$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 107
$ws.Range("A107").Value = 106
$ws.Range("B107").Value = "Tuesday, Jan 10"
$ws.Range("C107").Value = "4:00 AM"
$ws.Range("D107").Value = "QY5919"
$ws.Range("E107").Value = "Leipzig"
$ws.Range("F107").Value = "(LEJ)"
$ws.Range("G107").Value = "DHL (Delivered with Pride Livery) "
$ws.Range("H107").Value = "B752"
$ws.Range("I107").Value = "(D-ALEV)"
$ws.Range("J107").Value = "3:40 AM"
$ws.Range("L107").Value = "0 hours, -20 minutes"

# Row 108
$ws.Range("A108").Value = 107
$ws.Range("B108").Value = "Tuesday, Jan 10"
$ws.Range("C108").Value = "5:10 AM"
$ws.Range("D108").Value = "BO625"
$ws.Range("E108").Value = "Madrid"
$ws.Range("F108").Value = "(MAD)"
$ws.Range("G108").Value = "Bluebird Nordic "
$ws.Range("H108").Value = "B734"
$ws.Range("I108").Value = "(TF-BBO)"
$ws.Range("J108").Value = "5:12 AM"
$ws.Range("L108").Value = "0 hours, 2 minutes"

# Row 109
$ws.Range("A109").Value = 108
$ws.Range("B109").Value = "Tuesday, Jan 10"
$ws.Range("C109").Value = "6:15 AM"
$ws.Range("D109").Value = "W61001"
$ws.Range("E109").Value = "London"
$ws.Range("F109").Value = "(LTN)"
$ws.Range("G109").Value = "Wizz Air "
$ws.Range("H109").Value = "A21N"
$ws.Range("I109").Value = "(HA-LZJ)"
$ws.Range("J109").Value = "6:19 AM"
$ws.Range("L109").Value = "0 hours, 4 minutes"

# Row 110
$ws.Range("A110").Value = 109
$ws.Range("B110").Value = "Tuesday, Jan 10"
$ws.Range("C110").Value = "6:15 AM"
$ws.Range("D110").Value = "W61215"
$ws.Range("E110").Value = "Oslo"
$ws.Range("F110").Value = "(TRF)"
$ws.Range("G110").Value = "Wizz Air "
$ws.Range("H110").Value = "A321"
$ws.Range("I110").Value = "(HA-LXE)"
$ws.Range("J110").Value = "6:27 AM"
$ws.Range("L110").Value = "0 hours, 12 minutes"

# Row 111
$ws.Range("A111").Value = 110
$ws.Range("B111").Value = "Tuesday, Jan 10"
$ws.Range("C111").Value = "6:20 AM"
$ws.Range("D111").Value = "W61015"
$ws.Range("E111").Value = "Leeds"
$ws.Range("F111").Value = "(LBA)"
$ws.Range("G111").Value = "Wizz Air "
$ws.Range("H111").Value = "A321"
$ws.Range("I111").Value = "(HA-LXD)"
$ws.Range("J111").Value = "6:30 AM"
$ws.Range("L111").Value = "0 hours, 10 minutes"

# Row 112
$ws.Range("A112").Value = 111
$ws.Range("B112").Value = "Tuesday, Jan 10"
$ws.Range("C112").Value = "6:45 AM"
$ws.Range("D112").Value = "LH1363"
$ws.Range("E112").Value = "Frankfurt"
$ws.Range("F112").Value = "(FRA)"
$ws.Range("G112").Value = "Lufthansa "
$ws.Range("H112").Value = "CRJ9"
$ws.Range("I112").Value = "(D-ACNE)"
$ws.Range("J112").Value = "6:47 AM"
$ws.Range("L112").Value = "0 hours, 2 minutes"

# Row 113
$ws.Range("A113").Value = 112
$ws.Range("B113").Value = "Tuesday, Jan 10"
$ws.Range("C113").Value = "7:25 AM"
$ws.Range("D113").Value = "W61071"
$ws.Range("E113").Value = "Eindhoven"
$ws.Range("F113").Value = "(EIN)"
$ws.Range("G113").Value = "Wizz Air "
$ws.Range("H113").Value = "A321"
$ws.Range("I113").Value = "(HA-LTC)"
$ws.Range("J113").Value = "7:27 AM"
$ws.Range("L113").Value = "0 hours, 2 minutes"

# Row 114
$ws.Range("A114").Value = 113
$ws.Range("B114").Value = "Tuesday, Jan 10"
$ws.Range("C114").Value = "7:25 AM"
$ws.Range("D114").Value = "E44091"
$ws.Range("E114").Value = "Marsa Alam"
$ws.Range("F114").Value = "(RMF)"
$ws.Range("G114").Value = "Enter Air "
$ws.Range("H114").Value = "B738"
$ws.Range("I114").Value = "(SP-ENX)"
$ws.Range("J114").Value = "7:36 AM"
$ws.Range("L114").Value = "0 hours, 11 minutes"

# Row 115
$ws.Range("A115").Value = 114
$ws.Range("B115").Value = "Tuesday, Jan 10"
$ws.Range("C115").Value = "7:30 AM"
$ws.Range("D115").Value = "UNKNOWN"
$ws.Range("E115").Value = "Gran Canaria"
$ws.Range("F115").Value = "(LPA)"
$ws.Range("G115").Value = "Enter Air "
$ws.Range("H115").Value = "B738"
$ws.Range("I115").Value = "(SP-ENP)"
$ws.Range("J115").Value = "7:43 AM"
$ws.Range("L115").Value = "0 hours, 13 minutes"

# Row 116
$ws.Range("A116").Value = 115
$ws.Range("B116").Value = "Tuesday, Jan 10"
$ws.Range("C116").Value = "7:50 AM"
$ws.Range("D116").Value = "FR6391"
$ws.Range("E116").Value = "London"
$ws.Range("F116").Value = "(STN)"
$ws.Range("G116").Value = "Ryanair "
$ws.Range("H116").Value = "B738"
$ws.Range("I116").Value = "(SP-RSO)"
$ws.Range("J116").Value = "7:53 AM"
$ws.Range("L116").Value = "0 hours, 3 minutes"

# Row 117
$ws.Range("A117").Value = 116
$ws.Range("B117").Value = "Tuesday, Jan 10"
$ws.Range("C117").Value = "8:05 AM"
$ws.Range("D117").Value = "UNKNOWN"
$ws.Range("E117").Value = "Cologne"
$ws.Range("F117").Value = "(CGN)"
$ws.Range("G117").Value = "Enter Air "
$ws.Range("H117").Value = "B738"
$ws.Range("I117").Value = "(SP-ESB)"
$ws.Range("J117").Value = "7:56 AM"
$ws.Range("L117").Value = "0 hours, -9 minutes"

# Row 118
$ws.Range("A118").Value = 117
$ws.Range("B118").Value = "Tuesday, Jan 10"
$ws.Range("C118").Value = "8:30 AM"
$ws.Range("D118").Value = "3Z7108"
$ws.Range("E118").Value = "Marsa Alam"
$ws.Range("F118").Value = "(RMF)"
$ws.Range("G118").Value = "Smartwings "
$ws.Range("H118").Value = "B38M"
$ws.Range("I118").Value = "(OK-SWE)"
$ws.Range("J118").Value = "8:31 AM"
$ws.Range("L118").Value = "0 hours, 1 minutes"

# Row 119
$ws.Range("A119").Value = 118
$ws.Range("B119").Value = "Tuesday, Jan 10"
$ws.Range("C119").Value = "9:35 AM"
$ws.Range("D119").Value = "BO951"
$ws.Range("E119").Value = "Paris"
$ws.Range("F119").Value = "(CDG)"
$ws.Range("G119").Value = "Bluebird Nordic "
$ws.Range("H119").Value = "B734"
$ws.Range("I119").Value = "(TF-BBJ)"
$ws.Range("J119").Value = "9:26 AM"
$ws.Range("L119").Value = "0 hours, -9 minutes"

# Row 120
$ws.Range("A120").Value = 119
$ws.Range("B120").Value = "Tuesday, Jan 10"
$ws.Range("C120").Value = "10:30 AM"
$ws.Range("D120").Value = "UNKNOWN"
$ws.Range("E120").Value = "Budapest"
$ws.Range("F120").Value = "(BUD)"
$ws.Range("G120").Value = "Wizz Air "
$ws.Range("H120").Value = "A320"
$ws.Range("I120").Value = "(HA-LPO)"
$ws.Range("J120").Value = "11:14 AM"
$ws.Range("L120").Value = "0 hours, 44 minutes"

# Row 121
$ws.Range("A121").Value = 120
$ws.Range("B121").Value = "Tuesday, Jan 10"
$ws.Range("C121").Value = "11:15 AM"
$ws.Range("D121").Value = "FR7100"
$ws.Range("E121").Value = "Oslo"
$ws.Range("F121").Value = "(OSL)"
$ws.Range("G121").Value = "Ryanair "
$ws.Range("H121").Value = "B738"
$ws.Range("I121").Value = "(SP-RSN)"
$ws.Range("J121").Value = "11:27 AM"
$ws.Range("L121").Value = "0 hours, 12 minutes"

# Row 122
$ws.Range("A122").Value = 121
$ws.Range("B122").Value = "Tuesday, Jan 10"
$ws.Range("C122").Value = "11:25 AM"
$ws.Range("D122").Value = "FR6403"
$ws.Range("E122").Value = "Manchester"
$ws.Range("F122").Value = "(MAN)"
$ws.Range("G122").Value = "Ryanair "
$ws.Range("H122").Value = "B38M"
$ws.Range("I122").Value = "(EI-HGY)"
$ws.Range("J122").Value = "11:39 AM"
$ws.Range("L122").Value = "0 hours, 14 minutes"

# Row 123
$ws.Range("A123").Value = 122
$ws.Range("B123").Value = "Tuesday, Jan 10"
$ws.Range("C123").Value = "11:55 AM"
$ws.Range("D123").Value = "W61281"
$ws.Range("E123").Value = "Kutaisi"
$ws.Range("F123").Value = "(KUT)"
$ws.Range("G123").Value = "Wizz Air "
$ws.Range("H123").Value = "A321"
$ws.Range("I123").Value = "(HA-LXP)"
$ws.Range("J123").Value = "12:06 PM"
$ws.Range("L123").Value = "0 hours, 11 minutes"

# Row 124
$ws.Range("A124").Value = 123
$ws.Range("B124").Value = "Tuesday, Jan 10"
$ws.Range("C124").Value = "12:50 PM"
$ws.Range("D124").Value = "W61219"
$ws.Range("E124").Value = "Bergen"
$ws.Range("F124").Value = "(BGO)"
$ws.Range("G124").Value = "Wizz Air "
$ws.Range("H124").Value = "A321"
$ws.Range("I124").Value = "(HA-LXD)"
$ws.Range("J124").Value = "12:51 PM"
$ws.Range("L124").Value = "0 hours, 1 minutes"

# Row 125
$ws.Range("A125").Value = 124
$ws.Range("B125").Value = "Tuesday, Jan 10"
$ws.Range("C125").Value = "12:50 PM"
$ws.Range("D125").Value = "W61029"
$ws.Range("E125").Value = "Funchal"
$ws.Range("F125").Value = "(FNC)"
$ws.Range("G125").Value = "Wizz Air "
$ws.Range("H125").Value = "A21N"
$ws.Range("I125").Value = "(HA-LZJ)"
$ws.Range("J125").Value = "1:00 PM"
$ws.Range("L125").Value = "0 hours, 10 minutes"

# Row 126
$ws.Range("A126").Value = 125
$ws.Range("B126").Value = "Tuesday, Jan 10"
$ws.Range("C126").Value = "1:40 PM"
$ws.Range("D126").Value = "FR6389"
$ws.Range("E126").Value = "Dortmund"
$ws.Range("F126").Value = "(DTM)"
$ws.Range("G126").Value = "Ryanair "
$ws.Range("H126").Value = "B738"
$ws.Range("I126").Value = "(SP-RSO)"
$ws.Range("J126").Value = "2:15 PM"
$ws.Range("L126").Value = "0 hours, 35 minutes"

# Row 127
$ws.Range("A127").Value = 126
$ws.Range("B127").Value = "Tuesday, Jan 10"
$ws.Range("C127").Value = "2:15 PM"
$ws.Range("D127").Value = "LH1357"
$ws.Range("E127").Value = "Frankfurt"
$ws.Range("F127").Value = "(FRA)"
$ws.Range("G127").Value = "Lufthansa "
$ws.Range("H127").Value = "CRJ9"
$ws.Range("I127").Value = "(D-ACNT)"
$ws.Range("J127").Value = "2:34 PM"
$ws.Range("L127").Value = "0 hours, 19 minutes"

# Row 128
$ws.Range("A128").Value = 127
$ws.Range("B128").Value = "Tuesday, Jan 10"
$ws.Range("C128").Value = "3:00 PM"
$ws.Range("D128").Value = "LO3884"
$ws.Range("E128").Value = "Warsaw"
$ws.Range("F128").Value = "(WAW)"
$ws.Range("G128").Value = "LOT (Grzeski Livery) "
$ws.Range("H128").Value = "E195"
$ws.Range("I128").Value = "(SP-LNB)"
$ws.Range("J128").Value = "3:06 PM"
$ws.Range("L128").Value = "0 hours, 6 minutes"

# Row 129
$ws.Range("A129").Value = 128
$ws.Range("B129").Value = "Tuesday, Jan 10"
$ws.Range("C129").Value = "3:55 PM"
$ws.Range("D129").Value = "W61175"
$ws.Range("E129").Value = "Barcelona"
$ws.Range("F129").Value = "(BCN)"
$ws.Range("G129").Value = "Wizz Air "
$ws.Range("H129").Value = "A321"
$ws.Range("I129").Value = "(HA-LXE)"
$ws.Range("J129").Value = "3:58 PM"
$ws.Range("L129").Value = "0 hours, 3 minutes"

# Row 130
$ws.Range("A130").Value = 129
$ws.Range("B130").Value = "Tuesday, Jan 10"
$ws.Range("C130").Value = "4:45 PM"
$ws.Range("D130").Value = "FR2472"
$ws.Range("E130").Value = "London"
$ws.Range("F130").Value = "(STN)"
$ws.Range("G130").Value = "Ryanair "
$ws.Range("H130").Value = "B738"
$ws.Range("I130").Value = "(SP-RSN)"
$ws.Range("J130").Value = "4:52 PM"
$ws.Range("L130").Value = "0 hours, 7 minutes"

# Row 131
$ws.Range("A131").Value = 130
$ws.Range("B131").Value = "Tuesday, Jan 10"
$ws.Range("C131").Value = "4:50 PM"
$ws.Range("D131").Value = "KL1816"
$ws.Range("E131").Value = "Amsterdam"
$ws.Range("F131").Value = "(AMS)"
$ws.Range("G131").Value = "KLM "
$ws.Range("H131").Value = "E295"
$ws.Range("I131").Value = "(PH-NXN)"
$ws.Range("J131").Value = "4:59 PM"
$ws.Range("L131").Value = "0 hours, 9 minutes"

# Row 132
$ws.Range("A132").Value = 131
$ws.Range("B132").Value = "Tuesday, Jan 10"
$ws.Range("C132").Value = "6:20 PM"
$ws.Range("D132").Value = "LO3886"
$ws.Range("E132").Value = "Warsaw"
$ws.Range("F132").Value = "(WAW)"
$ws.Range("G132").Value = "LOT "
$ws.Range("H132").Value = "E190"
$ws.Range("I132").Value = "(SP-LMG)"
$ws.Range("J132").Value = "6:16 PM"
$ws.Range("L132").Value = "0 hours, -4 minutes"

# Row 133
$ws.Range("A133").Value = 132
$ws.Range("B133").Value = "Tuesday, Jan 10"
$ws.Range("C133").Value = "6:35 PM"
$ws.Range("D133").Value = "FR823"
$ws.Range("E133").Value = "Venice"
$ws.Range("F133").Value = "(VCE)"
$ws.Range("G133").Value = "Ryanair "
$ws.Range("H133").Value = "B738"
$ws.Range("I133").Value = "(9H-QBP)"
$ws.Range("J133").Value = "6:43 PM"
$ws.Range("L133").Value = "0 hours, 8 minutes"

# Row 134
$ws.Range("A134").Value = 133
$ws.Range("B134").Value = "Tuesday, Jan 10"
$ws.Range("C134").Value = "9:35 PM"
$ws.Range("D134").Value = "SAR1981"
$ws.Range("E134").Value = "Cologne"
$ws.Range("F134").Value = "(CGN)"
$ws.Range("G134").Value = "SprintAir "
$ws.Range("H134").Value = "AT75"
$ws.Range("I134").Value = "(SP-SPF)"
$ws.Range("J134").Value = "9:42 PM"
$ws.Range("L134").Value = "0 hours, 7 minutes"

# Row 135
$ws.Range("A135").Value = 134
$ws.Range("B135").Value = "Tuesday, Jan 10"
$ws.Range("C135").Value = "9:52 PM"
$ws.Range("D135").Value = "3V4563"
$ws.Range("E135").Value = "Paris"
$ws.Range("F135").Value = "(CDG)"
$ws.Range("G135").Value = "FedEx "
$ws.Range("H135").Value = "B738"
$ws.Range("I135").Value = "(OE-IWF)"
$ws.Range("J135").Value = "9:57 PM"
$ws.Range("L135").Value = "0 hours, 5 minutes"

# Row 136
$ws.Range("A136").Value = 135
$ws.Range("B136").Value = "Tuesday, Jan 10"
$ws.Range("C136").Value = "10:20 PM"
$ws.Range("D136").Value = "QY5917"
$ws.Range("E136").Value = "Leipzig"
$ws.Range("F136").Value = "(LEJ)"
$ws.Range("G136").Value = "ASL Airlines "
$ws.Range("H136").Value = "B734"
$ws.Range("I136").Value = "(EI-STU)"
$ws.Range("J136").Value = "10:26 PM"
$ws.Range("L136").Value = "0 hours, 6 minutes"
